$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (A15:F15) already holds the exact same id/nama/kota/handphone/tamu
# values as the new guest we need to append in row 23 ("DDMMYYFN20",
# "Fahrezi Rizqiawan", "Kota Bekasi", "089662690020", 5). Copy it down so
# the text-vs-number typing (e.g. the leading zero in the phone number)
# round-trips exactly, instead of retyping values that Excel would
# otherwise re-infer as numbers.
$ws.Range("A15:F15").Copy($ws.Range("A23:F23"))

# Only the "hubungan" (relationship) column differs for this guest.
$ws.Range("F23").Value = "Kerabat Mempelai Pria"
